$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("accuracy")
$ws1.Range("B2").Value = [double]"1.189224106458005e-08"
$ws1.Range("B3").Value = [double]"0.743532461668439"
$ws1.Range("B4").Value = [double]"0.2853081108123064"
$ws1.Range("B5").Value = [double]"0.3835110692831069"
$ws1.Range("B6").Value = [double]"0.4423783534346787"
$ws1.Range("B7").Value = [double]"0.01176300432789639"
$ws1.Range("B8").Value = [double]"0.1183057268373853"
$ws1.Range("B9").Value = [double]"0.5770324017264793"
$ws1.Range("B10").Value = [double]"2.143314765724455e-05"
$ws1.Range("B11").Value = [double]"3.863341453943772e-05"
$ws1.Range("B12").Value = [double]"3.44810298557851e-05"
$ws1.Range("B13").Value = [double]"0.02182190128758563"
$ws1.Range("B14").Value = [double]"0.001254375329880237"
$ws1.Range("B15").Value = [double]"1.553525709426425e-08"

$ws2 = $wb.Worksheets.Item("sensitivity")
$ws2.Range("B2").Value = [double]"2.100648757946631e-05"
$ws2.Range("B3").Value = [double]"8.056838827842495e-09"
$ws2.Range("B4").Value = [double]"5.433299535704519e-07"
$ws2.Range("B5").Value = [double]"0.007032931722704711"
$ws2.Range("B6").Value = [double]"0.00746958453953499"
$ws2.Range("B7").Value = [double]"7.061502350243177e-12"
$ws2.Range("B8").Value = [double]"2.358169682807503e-07"
$ws2.Range("B9").Value = [double]"1.002723295488438e-08"
$ws2.Range("B10").Value = [double]"0.002949310675381848"
$ws2.Range("B11").Value = [double]"3.555597563440914e-08"
$ws2.Range("B12").Value = [double]"4.066554996949832e-14"
$ws2.Range("B13").Value = [double]"6.049066718272365e-10"
$ws2.Range("B14").Value = [double]"4.955163747957975e-13"
$ws2.Range("B15").Value = [double]"1.16992617465201e-16"

$ws3 = $wb.Worksheets.Item("specificity")
$ws3.Range("B2").Value = [double]"2.388614298683008e-07"
$ws3.Range("B3").Value = [double]"1.748864560469862e-09"
$ws3.Range("B4").Value = [double]"3.421516618467542e-10"
$ws3.Range("B5").Value = [double]"0.0006911052184860778"
$ws3.Range("B6").Value = [double]"0.0002927074424314717"
$ws3.Range("B7").Value = [double]"7.966900069562025e-08"
$ws3.Range("B8").Value = [double]"0.0002574540638237679"
$ws3.Range("B9").Value = [double]"2.10491016137352e-06"
$ws3.Range("B10").Value = [double]"9.401525086603872e-06"
$ws3.Range("B11").Value = [double]"0.003634820685442738"
$ws3.Range("B12").Value = [double]"1.00709518339693e-12"
$ws3.Range("B13").Value = [double]"2.303678704473348e-08"
$ws3.Range("B14").Value = [double]"1.328322046511736e-09"
$ws3.Range("B15").Value = [double]"4.580073111445163e-14"

$ws4 = $wb.Worksheets.Item("time")
$timeVal = [double]"3.54161358920777e-34"
for ($r = 2; $r -le 15; $r++) {
    $ws4.Cells.Item($r, 2).Value = $timeVal
}
